# Update "Förändrad" date column (C) for all data rows (2-46) from
# serial 45202 (2023-10-03) to serial 45203 (2023-10-04).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C46")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
